$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# Seed the new rows (1060-1080) by cloning the format of the last existing row (1059),
# then overwrite the values - this keeps the 's="1" ("import")' cell style used
# throughout the sheet instead of Excel inventing a new style index.
$ws.Range("A1059:C1059").Copy() | Out-Null
$ws.Range("A1060:C1080").PasteSpecial(-4122) | Out-Null

$ws.Range("A1060").Value = "cs"
$ws.Range("B1060").Value = 'lab.vendor.label'
$ws.Range("C1060").Value = 'Výrobci'

$ws.Range("A1061").Value = "cs"
$ws.Range("B1061").Value = 'lab.vendor.title'
$ws.Range("C1061").Value = 'Výrobci'

$ws.Range("A1062").Value = "cs"
$ws.Range("B1062").Value = 'lab.vendor.table.name'
$ws.Range("C1062").Value = 'Jméno'

$ws.Range("A1063").Value = "cs"
$ws.Range("B1063").Value = 'lab.vendor.button.create'
$ws.Range("C1063").Value = 'Nový výrobce'

$ws.Range("A1064").Value = "cs"
$ws.Range("B1064").Value = 'lab.vendor.filter.title'
$ws.Range("C1064").Value = 'Filtrovat výrobce'

$ws.Range("A1065").Value = "cs"
$ws.Range("B1065").Value = 'lab.vendor.button.create'
$ws.Range("C1065").Value = 'Nový výrobce'

$ws.Range("A1066").Value = "cs"
$ws.Range("B1066").Value = 'lab.vendor.context.menu'
$ws.Range("C1066").Value = 'Výrobce [{{data.name}}]'

$ws.Range("A1067").Value = "cs"
$ws.Range("B1067").Value = 'lab.vendor.preview'
$ws.Range("C1067").Value = 'Náhled výrobce'

$ws.Range("A1068").Value = "cs"
$ws.Range("B1068").Value = 'lab.vendor.button.edit'
$ws.Range("C1068").Value = 'Upravit výrobce'

$ws.Range("A1069").Value = "cs"
$ws.Range("B1069").Value = 'lab.vendor.button.delete'
$ws.Range("C1069").Value = 'Odstranit výrobce'

$ws.Range("A1070").Value = "cs"
$ws.Range("B1070").Value = 'lab.vendor.button.delete.confirm.title'
$ws.Range("C1070").Value = 'Odstranit výrobce'

$ws.Range("A1071").Value = "cs"
$ws.Range("B1071").Value = 'lab.vendor.button.delete.confirm.ok'
$ws.Range("C1071").Value = 'Odstranit výrobce'

$ws.Range("A1072").Value = "cs"
$ws.Range("B1072").Value = 'lab.vendor.button.delete.confirm'
$ws.Range("C1072").Value = 'Opravdu si přejete odstranit vybraného výrobce? Tento krok s sebou efektivně může vzít opravdu velké množství dat ze systému, jelikož výrobci jsou použiti u atomizérů, modů, drátů, zkrátka všude. Díky tomu může dojít k brutálnímu promazání dat. Před tím, než budete pokračovat, se hezky pomodlete ke svému oblíbenému božstvu.'

$ws.Range("A1073").Value = "cs"
$ws.Range("B1073").Value = 'lab.vendor.table.footer.label'
$ws.Range("C1073").Value = 'Počet výrobců [{{data.total}}]'

$ws.Range("A1074").Value = "cs"
$ws.Range("B1074").Value = 'lab.vendor.preview.name'
$ws.Range("C1074").Value = 'Jméno'

$ws.Range("A1075").Value = "cs"
$ws.Range("B1075").Value = 'lab.vendor.common.tab'
$ws.Range("C1075").Value = 'Výrobce'

$ws.Range("A1076").Value = "cs"
$ws.Range("B1076").Value = 'lab.vendor.atomizers.tab'
$ws.Range("C1076").Value = 'Atomizéry'

$ws.Range("A1077").Value = "cs"
$ws.Range("B1077").Value = 'lab.vendor.mods.tab'
$ws.Range("C1077").Value = 'Mody'

$ws.Range("A1078").Value = "cs"
$ws.Range("B1078").Value = 'lab.vendor.liquid.tab'
$ws.Range("C1078").Value = 'Liquidy'

$ws.Range("A1079").Value = "cs"
$ws.Range("B1079").Value = 'lab.vendor.index.label'
$ws.Range("C1079").Value = 'Výrobce [{{data.name}}]'

$ws.Range("A1080").Value = "cs"
$ws.Range("B1080").Value = 'lab.vendor.index.title'
$ws.Range("C1080").Value = 'Detail výrobce'

# Row 1072 carries the long "delete confirm" message and wraps to a taller row
# (mirrors the analogous "wire" confirm row which is ht="30" for a shorter string).
$ws.Range("A1072:C1072").EntireRow.RowHeight = 60

# Restore the sheet view to the area the author was looking at after the edit.
$ws.Activate()
$ws.Range("B1071").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1053
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 100
